# "Done with section 4 need to make changes in comments"
#
# Section-4 style placeholder cells (bold, empty) get extended out to the
# right on the header/total rows (5 & 6), and the "Section 4 comment"
# answers that were mistakenly recorded one row too early (rows 10 & 13)
# get moved down to their correct rows (12 & 15 respectively).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 5 & 6: extend the bold/empty formatting block further right.
$ws.Range("AX5:BR5").Font.Bold = $true
$ws.Range("AX6:BQ6").Font.Bold = $true

# Row 10 -> Row 12: the "section 4 comments" answers belong on row 12, not 10.
$src10 = $ws.Range("AX10:BR10")
$dst12 = $ws.Range("AX12:BR12")
$dst12.Value = $src10.Value()
$src10.ClearContents()

# Row 13 -> Row 15: same fix, the answers belong on row 15, not 13.
$src13 = $ws.Range("AX13:BP13")
$dst15 = $ws.Range("AX15:BP15")
$dst15.Value = $src13.Value()
$src13.ClearContents()
